$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 15.33839062041541
$ws.Cells.Item(2, 3).Value = 9.166594144469427
$ws.Cells.Item(2, 4).Value = 14.07034218581211
$ws.Cells.Item(2, 5).Value = 14.79430189541348
$ws.Cells.Item(2, 7).Value = 41.44253356168132
$ws.Cells.Item(2, 8).Value = 17.58203912378771
$ws.Cells.Item(2, 10).Value = 8.826199713099319
$ws.Cells.Item(2, 11).Value = 11.10866066631313
$ws.Cells.Item(2, 12).Value = 11.8215222220309
$ws.Cells.Item(2, 14).Value = 20.93027811675206
$ws.Cells.Item(2, 15).Value = 28.42691061579772
$ws.Cells.Item(3, 2).Value = 15.15663643075777
$ws.Cells.Item(3, 3).Value = 9.165480868125467
$ws.Cells.Item(3, 4).Value = 14.06382243465257
$ws.Cells.Item(3, 5).Value = 14.81636582425137
$ws.Cells.Item(3, 7).Value = 41.53610834912766
$ws.Cells.Item(3, 8).Value = 17.62860842368504
$ws.Cells.Item(3, 10).Value = 8.838320442575942
$ws.Cells.Item(3, 11).Value = 10.97476941970377
$ws.Cells.Item(3, 12).Value = 11.82261383206795
$ws.Cells.Item(3, 14).Value = 20.98957976191751
$ws.Cells.Item(3, 15).Value = 28.50369667394247
$ws.Cells.Item(4, 2).Value = 15.04668084953194
$ws.Cells.Item(4, 3).Value = 9.165083734739749
$ws.Cells.Item(4, 4).Value = 14.06231482897556
$ws.Cells.Item(4, 5).Value = 14.83191479966953
$ws.Cells.Item(4, 7).Value = 41.60278602852386
$ws.Cells.Item(4, 8).Value = 17.65951432980759
$ws.Cells.Item(4, 10).Value = 8.846186883909763
$ws.Cells.Item(4, 11).Value = 10.89344953719039
$ws.Cells.Item(4, 12).Value = 11.82480708734968
$ws.Cells.Item(4, 14).Value = 21.02770541563973
$ws.Cells.Item(4, 15).Value = 28.55555840136425
$ws.Cells.Item(5, 2).Value = 15.00233554053788
$ws.Cells.Item(5, 3).Value = 9.164994345513232
$ws.Cells.Item(5, 4).Value = 14.06233001505115
$ws.Cells.Item(5, 5).Value = 14.83875490236781
$ws.Cells.Item(5, 7).Value = 41.63227234944149
$ws.Cells.Item(5, 8).Value = 17.67269049347093
$ws.Cells.Item(5, 10).Value = 8.849499510513446
$ws.Cells.Item(5, 11).Value = 10.86057038679627
$ws.Cells.Item(5, 12).Value = 11.82608498272208
$ws.Cells.Item(5, 14).Value = 21.04367423421088
$ws.Cells.Item(5, 15).Value = 28.57787739107056
$ws.Cells.Item(6, 2).Value = 14.99500140112885
$ws.Cells.Item(6, 3).Value = 9.164983890356995
$ws.Cells.Item(6, 4).Value = 14.06237061269991
$ws.Cells.Item(6, 5).Value = 14.83992113449052
$ws.Cells.Item(6, 7).Value = 41.6373081897915
$ws.Cells.Item(6, 8).Value = 17.67491353208473
$ws.Cells.Item(6, 10).Value = 8.850056040369751
$ws.Cells.Item(6, 11).Value = 10.85512752643016
$ws.Cells.Item(6, 12).Value = 11.82632040820547
$ws.Cells.Item(6, 14).Value = 21.04635199145783
$ws.Cells.Item(6, 15).Value = 28.58165498807267
$ws.Cells.Item(7, 2).Value = 15.04608085573448
$ws.Cells.Item(7, 3).Value = 9.16508223531593
$ws.Cells.Item(7, 4).Value = 14.06231248232276
$ws.Cells.Item(7, 5).Value = 14.83200500739462
$ws.Cells.Item(7, 7).Value = 41.60317432573544
$ws.Cells.Item(7, 8).Value = 17.65968967220154
$ws.Cells.Item(7, 10).Value = 8.846231125555894
$ws.Cells.Item(7, 11).Value = 10.89300502058179
$ws.Cells.Item(7, 12).Value = 11.82482276473162
$ws.Cells.Item(7, 14).Value = 21.02791902466179
$ws.Cells.Item(7, 15).Value = 28.55585460630247
$ws.Cells.Item(8, 2).Value = 15.27541085573456
$ws.Cells.Item(8, 3).Value = 9.166151092922313
$ws.Cells.Item(8, 4).Value = 14.06757752831311
$ws.Cells.Item(8, 5).Value = 14.80149437875584
$ws.Cells.Item(8, 7).Value = 41.47288152470301
$ws.Cells.Item(8, 8).Value = 17.59761653572248
$ws.Cells.Item(8, 10).Value = 8.83029106714835
$ws.Cells.Item(8, 11).Value = 11.06233172076797
$ws.Cells.Item(8, 12).Value = 11.82158323179333
$ws.Cells.Item(8, 14).Value = 20.95037029969135
$ws.Cells.Item(8, 15).Value = 28.45240751999965
$ws.Cells.Item(9, 2).Value = 15.73595354289211
$ws.Cells.Item(9, 3).Value = 9.170501623265402
$ws.Cells.Item(9, 4).Value = 14.09760243227134
$ws.Cells.Item(9, 5).Value = 14.75752456394244
$ws.Cells.Item(9, 7).Value = 41.29073300237437
$ws.Cells.Item(9, 8).Value = 17.49422199037434
$ws.Cells.Item(9, 10).Value = 8.802385237399479
$ws.Cells.Item(9, 11).Value = 11.39986687117076
$ws.Cells.Item(9, 12).Value = 11.82726098057809
$ws.Cells.Item(9, 14).Value = 20.81184068804342
$ws.Cells.Item(9, 15).Value = 28.2869856713092
$ws.Cells.Item(10, 2).Value = 16.07795034702833
$ws.Cells.Item(10, 3).Value = 9.175048627501981
$ws.Cells.Item(10, 4).Value = 14.13151845647269
$ws.Cells.Item(10, 5).Value = 14.73485987532015
$ws.Cells.Item(10, 7).Value = 41.2018468123768
$ws.Cells.Item(10, 8).Value = 17.42941185486354
$ws.Cells.Item(10, 10).Value = 8.783907455975818
$ws.Cells.Item(10, 11).Value = 11.64908984749292
$ws.Cells.Item(10, 12).Value = 11.83868952418934
$ws.Cells.Item(10, 14).Value = 20.71823636293386
$ws.Cells.Item(10, 15).Value = 28.18830927972106
$ws.Cells.Item(11, 2).Value = 16.2336981634975
$ws.Cells.Item(11, 3).Value = 9.177405378295663
$ws.Cells.Item(11, 4).Value = 14.14948303306288
$ws.Cells.Item(11, 5).Value = 14.72663488412727
$ws.Cells.Item(11, 7).Value = 41.17120394327674
$ws.Cells.Item(11, 8).Value = 17.40234553724923
$ws.Cells.Item(11, 10).Value = 8.775936995372378
$ws.Cells.Item(11, 11).Value = 11.76229438917656
$ws.Cells.Item(11, 12).Value = 11.84544583862034
$ws.Cells.Item(11, 14).Value = 20.67741040862082
$ws.Cells.Item(11, 15).Value = 28.14838904231184
$ws.Cells.Item(12, 2).Value = 16.29264914894304
$ws.Cells.Item(12, 3).Value = 9.178338825522273
$ws.Cells.Item(12, 4).Value = 14.15664653089336
$ws.Cells.Item(12, 5).Value = 14.72381940820119
$ws.Cells.Item(12, 7).Value = 41.16101024985332
$ws.Cells.Item(12, 8).Value = 17.39244333755358
$ws.Cells.Item(12, 10).Value = 8.772981058485962
$ws.Cells.Item(12, 11).Value = 11.80510164565433
$ws.Cells.Item(12, 12).Value = 11.84822643463293
$ws.Cells.Item(12, 14).Value = 20.66220176179375
$ws.Cells.Item(12, 15).Value = 28.13398715221317
$ws.Cells.Item(13, 2).Value = 16.27995503173681
$ws.Cells.Item(13, 3).Value = 9.178135974984373
$ws.Cells.Item(13, 4).Value = 14.15508776784757
$ws.Cells.Item(13, 5).Value = 14.72441247863466
$ws.Cells.Item(13, 7).Value = 41.16314290408099
$ws.Cells.Item(13, 8).Value = 17.39456051379217
$ws.Cells.Item(13, 10).Value = 8.773614905752765
$ws.Cells.Item(13, 11).Value = 11.79588561751614
$ws.Cells.Item(13, 12).Value = 11.84761773475963
$ws.Cells.Item(13, 14).Value = 20.66546605918067
$ws.Cells.Item(13, 15).Value = 28.13705704637584
$ws.Cells.Item(14, 2).Value = 16.23854893844457
$ws.Cells.Item(14, 3).Value = 9.177481354572494
$ws.Cells.Item(14, 4).Value = 14.15006517323814
$ws.Cells.Item(14, 5).Value = 14.72639726214824
$ws.Cells.Item(14, 7).Value = 41.17033703216168
$ws.Cells.Item(14, 8).Value = 17.40152391958629
$ws.Cells.Item(14, 10).Value = 8.775692561570377
$ws.Cells.Item(14, 11).Value = 11.76581759423836
$ws.Cells.Item(14, 12).Value = 11.84567015907271
$ws.Cells.Item(14, 14).Value = 20.67615415497397
$ws.Cells.Item(14, 15).Value = 28.14718985644554
$ws.Cells.Item(15, 2).Value = 16.21318142564578
$ws.Cells.Item(15, 3).Value = 9.177085705947208
$ws.Cells.Item(15, 4).Value = 14.14703553966449
$ws.Cells.Item(15, 5).Value = 14.72765193516463
$ws.Cells.Item(15, 7).Value = 41.17492732902461
$ws.Cells.Item(15, 8).Value = 17.4058344196423
$ws.Cells.Item(15, 10).Value = 8.776973291422351
$ws.Cells.Item(15, 11).Value = 11.74739106467973
$ws.Cells.Item(15, 12).Value = 11.84450608388114
$ws.Cells.Item(15, 14).Value = 20.68273360990322
$ws.Cells.Item(15, 15).Value = 28.15348963130794
$ws.Cells.Item(16, 2).Value = 16.06777084526971
$ws.Cells.Item(16, 3).Value = 9.174900370348176
$ws.Cells.Item(16, 4).Value = 14.13039512803393
$ws.Cells.Item(16, 5).Value = 14.73543929987556
$ws.Cells.Item(16, 7).Value = 41.20404662213618
$ws.Cells.Item(16, 8).Value = 17.4312293033553
$ws.Cells.Item(16, 10).Value = 8.784437077667423
$ws.Cells.Item(16, 11).Value = 11.64168512655431
$ws.Cells.Item(16, 12).Value = 11.83827916874333
$ws.Cells.Item(16, 14).Value = 20.72093965590612
$ws.Cells.Item(16, 15).Value = 28.19101818811182
$ws.Cells.Item(17, 2).Value = 15.97857280923133
$ws.Cells.Item(17, 3).Value = 9.17363325154284
$ws.Cells.Item(17, 4).Value = 14.12083358497746
$ws.Cells.Item(17, 5).Value = 14.74075020355242
$ws.Cells.Item(17, 7).Value = 41.22441995146435
$ws.Cells.Item(17, 8).Value = 17.44742692443338
$ws.Cells.Item(17, 10).Value = 8.789127131699438
$ws.Cells.Item(17, 11).Value = 11.57676852681507
$ws.Cells.Item(17, 12).Value = 11.83485678593367
$ws.Cells.Item(17, 14).Value = 20.74482655031849
$ws.Cells.Item(17, 15).Value = 28.2153136731986
$ws.Cells.Item(18, 2).Value = 15.92728767915028
$ws.Cells.Item(18, 3).Value = 9.172931601868866
$ws.Cells.Item(18, 4).Value = 14.11557299800217
$ws.Cells.Item(18, 5).Value = 14.74400121914027
$ws.Cells.Item(18, 7).Value = 41.23705972971423
$ws.Cells.Item(18, 8).Value = 17.45697079879606
$ws.Cells.Item(18, 10).Value = 8.79186570367367
$ws.Cells.Item(18, 11).Value = 11.53941652521639
$ws.Cells.Item(18, 12).Value = 11.83303502977759
$ws.Cells.Item(18, 14).Value = 20.75873092741091
$ws.Cells.Item(18, 15).Value = 28.22975538686047
$ws.Cells.Item(19, 2).Value = 15.90992831894192
$ws.Cells.Item(19, 3).Value = 9.172698714642488
$ws.Cells.Item(19, 4).Value = 14.11383301084192
$ws.Cells.Item(19, 5).Value = 14.74513569661197
$ws.Cells.Item(19, 7).Value = 41.24149755121792
$ws.Cells.Item(19, 8).Value = 17.46024125854201
$ws.Cells.Item(19, 10).Value = 8.79279998330991
$ws.Cells.Item(19, 11).Value = 11.52676855838879
$ws.Cells.Item(19, 12).Value = 11.83244346541939
$ws.Cells.Item(19, 14).Value = 20.76346712995763
$ws.Cells.Item(19, 15).Value = 28.2347253922541
$ws.Cells.Item(20, 2).Value = 15.98806648340847
$ws.Cells.Item(20, 3).Value = 9.173765330494531
$ws.Cells.Item(20, 4).Value = 14.1218267238654
$ws.Cells.Item(20, 5).Value = 14.74016453491095
$ws.Cells.Item(20, 7).Value = 41.2221557791766
$ws.Cells.Item(20, 8).Value = 17.44567912325545
$ws.Cells.Item(20, 10).Value = 8.788623628475287
$ws.Cells.Item(20, 11).Value = 11.58368069478304
$ws.Cells.Item(20, 12).Value = 11.83520593357436
$ws.Cells.Item(20, 14).Value = 20.74226665329546
$ws.Cells.Item(20, 15).Value = 28.21267897840828
$ws.Cells.Item(21, 2).Value = 16.25071206273666
$ws.Cells.Item(21, 3).Value = 9.177672523499375
$ws.Cells.Item(21, 4).Value = 14.15153067493798
$ws.Cells.Item(21, 5).Value = 14.72580617030125
$ws.Cells.Item(21, 7).Value = 41.16818565968786
$ws.Cells.Item(21, 8).Value = 17.39946917655391
$ws.Cells.Item(21, 10).Value = 8.775080614964391
$ws.Cells.Item(21, 11).Value = 11.77465124068571
$ws.Cells.Item(21, 12).Value = 11.84623619596746
$ws.Cells.Item(21, 14).Value = 20.67300799207987
$ws.Cells.Item(21, 15).Value = 28.14419419402453
$ws.Cells.Item(22, 2).Value = 16.42218350927644
$ws.Cells.Item(22, 3).Value = 9.180464844489093
$ws.Cells.Item(22, 4).Value = 14.17304430269892
$ws.Cells.Item(22, 5).Value = 14.71816544298082
$ws.Cells.Item(22, 7).Value = 41.14113267953427
$ws.Cells.Item(22, 8).Value = 17.37129211141856
$ws.Cells.Item(22, 10).Value = 8.766592500195207
$ws.Cells.Item(22, 11).Value = 11.89909051663517
$ws.Cells.Item(22, 12).Value = 11.85473890153565
$ws.Cells.Item(22, 14).Value = 20.62920744966369
$ws.Cells.Item(22, 15).Value = 28.10360344540277
$ws.Cells.Item(23, 2).Value = 16.33069957476836
$ws.Cells.Item(23, 3).Value = 9.178952838148083
$ws.Cells.Item(23, 4).Value = 14.16137128629449
$ws.Cells.Item(23, 5).Value = 14.72208418470722
$ws.Cells.Item(23, 7).Value = 41.15481876150967
$ws.Cells.Item(23, 8).Value = 17.38614563350499
$ws.Cells.Item(23, 10).Value = 8.77108963965031
$ws.Cells.Item(23, 11).Value = 11.83272075792556
$ws.Cells.Item(23, 12).Value = 11.85008309981434
$ws.Cells.Item(23, 14).Value = 20.65245103778388
$ws.Cells.Item(23, 15).Value = 28.12488593299777
$ws.Cells.Item(24, 2).Value = 15.98377440175826
$ws.Cells.Item(24, 3).Value = 9.173705533985792
$ws.Cells.Item(24, 4).Value = 14.1213769889045
$ws.Cells.Item(24, 5).Value = 14.74042869982295
$ws.Cells.Item(24, 7).Value = 41.22317652445732
$ws.Cells.Item(24, 8).Value = 17.44646858258235
$ws.Cells.Item(24, 10).Value = 8.788851130821877
$ws.Cells.Item(24, 11).Value = 11.58055579609774
$ws.Cells.Item(24, 12).Value = 11.83504762958759
$ws.Cells.Item(24, 14).Value = 20.74342344848439
$ws.Cells.Item(24, 15).Value = 28.21386864774728
$ws.Cells.Item(25, 2).Value = 15.61051189918863
$ws.Cells.Item(25, 3).Value = 9.169085763843601
$ws.Cells.Item(25, 4).Value = 14.08738534461943
$ws.Cells.Item(25, 5).Value = 14.76772429991967
$ws.Cells.Item(25, 7).Value = 41.33213102920185
$ws.Cells.Item(25, 8).Value = 17.52023269670111
$ws.Cells.Item(25, 10).Value = 8.809577580255011
$ws.Cells.Item(25, 11).Value = 11.30818586953269
$ws.Cells.Item(25, 12).Value = 11.82444454546773
$ws.Cells.Item(25, 14).Value = 20.68273360990322
$ws.Cells.Item(25, 15).Value = 28.15348963130794
